# completed the train and tests
# Fill in the new Train class test-plan rows (7-13), columns E (Preconditions),
# F (Method Inputs) and G (Expected Result) of the Table1 test plan.
# Cells are populated in G, F, E order per row to reproduce the same
# shared-string insertion order as the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - __init__ happy path
$ws.Range("G7").Value = "The train instance is created successfully with the attributes correctly set."
$ws.Range("F7").Value = 'make = "Siemens",                 model = "Intercity Subway",       cars = 13                           base_fuel_rate = 0.03'
$ws.Range("E7").Value = "None"

# Row 8 - __init__ blank make
$ws.Range("G8").Value = "ValueError"
$ws.Range("F8").Value = 'make = "  ",                                model = "Intercity Subway",       cars = 13                           base_fuel_rate = 0.03'
$ws.Range("E8").Value = "None"

# Row 9 - __init__ blank model
$ws.Range("G9").Value = "ValueError"
$ws.Range("F9").Value = 'make = "Siemens",                 model = "    ",                            cars = 13                           base_fuel_rate = 0.03'
$ws.Range("E9").Value = "None"

# Row 10 - __init__ cars not an integer
$ws.Range("G10").Value = "ValueError"
$ws.Range("F10").Value = 'make = "Siemens",                 model = "Intercity Subway",       cars = "cars"                       base_fuel_rate = 0.03'
$ws.Range("E10").Value = "None"

# Row 11 - __init__ base_fuel_rate not numeric
$ws.Range("G11").Value = "ValueError"
$ws.Range("F11").Value = 'make = "Siemens",                 model = "Intercity Subway",       cars = 13                           base_fuel_rate = None'
$ws.Range("E11").Value = "None"

# Row 12 - __str__ returns the formatted string (leading apostrophe forces text /
# reproduces the quotePrefix formatting applied to this cell)
$ws.Range("G12").Formula = '''Make: Siemens \n Model: Intercity Subway\nThis train has a base fuel rate of 0.03 litres/kilometer, and contains 13 cars.'''
$ws.Range("F12").Value = "None"
$ws.Range("E12").Value = 'make = "Siemens",                 model = "Intercity Subway",       cars = 13                           base_fuel_rate = 0.03'

# Row 13 - calculate_fuel_requirements returns correct calculated value
$ws.Range("G13").Value = 42.9
$ws.Range("F13").Value = "distance = 100.0"
$ws.Range("E13").Value = 'make = "Siemens",                 model = "Intercity Subway",       cars = 13                           base_fuel_rate = 0.03'

# Narrow columns E and F now that the real (shorter) content is in place
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666
$ws.Columns.Item(6).ColumnWidth = 22.498697916666668

# Match the author's final selection
$ws.Range("E12").Select()
